$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-less approach: set date cells to text format first so Excel does not
# auto-convert the dash-separated dates into date serial numbers, write the
# new text, then restore the "Normal" style so the cell's style index stays
# the same as before (no explicit style, same as the original inlineStr cells).
function Set-TextDate($addr, $text) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextDate "A3" "28-07-2022"
Set-TextDate "A4" "01-08-2022"
Set-TextDate "A5" "04-08-2022"
Set-TextDate "A6" "08-08-2022"
Set-TextDate "A7" "11-08-2022"
Set-TextDate "A8" "15-08-2022"
Set-TextDate "A9" "18-08-2022"
Set-TextDate "A10" "22-08-2022"
Set-TextDate "A11" "25-08-2022"
Set-TextDate "A12" "29-08-2022"
Set-TextDate "A13" "01-09-2022"
Set-TextDate "A14" "05-09-2022"
Set-TextDate "A15" "08-09-2022"
Set-TextDate "A16" "12-09-2022"
Set-TextDate "A17" "15-09-2022"
Set-TextDate "A18" "19-09-2022"
Set-TextDate "A19" "22-09-2022"
Set-TextDate "A20" "26-09-2022"
Set-TextDate "A21" "29-09-2022"

# Update attendance counts for row 3 (28-07-2022): D3 0->1, G3 0->1
$ws.Range("D3").Value = 1
$ws.Range("G3").Value = 1

# Update attendance counts for row 4 (01-08-2022): D4 0->1, E4 0->1, H4 1->0
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("H4").Value = 0
